$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cryptos list refresh: update Price (D) and Volume(1h) (E) columns per row.
# For D-column values that look like plain numbers (single decimal point),
# force text format first so Excel keeps them as literal strings (e.g. "45.60"
# keeps its trailing zero) instead of silently converting to a numeric value.

$ws.Range("D2").Value = "43.797.53"
$ws.Range("E2").Value = "  +0.39%  "
$ws.Range("D3").Value = "2.340.20"
$ws.Range("E3").Value = "  +4.52%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "96.56"
$ws.Range("E5").Value = "  +3.23%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "271.84"
$ws.Range("E6").Value = "  +0.66%  "
$ws.Range("E7").Value = "  +0.52%  "
$ws.Range("E8").Value = "  +0.04%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.625"
$ws.Range("E9").Value = "  +0.23%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "45.60"
$ws.Range("E10").Value = "  -1.92%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0942"
$ws.Range("E11").Value = "  +1.92%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "8.11"
$ws.Range("E12").Value = "  -1.05%  "
$ws.Range("E13").Value = "  +0.41%  "
$ws.Range("D14").Value = "2.686.38"
$ws.Range("E14").Value = "  +4.42%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.67"
$ws.Range("E15").Value = "  +3.37%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.870"
$ws.Range("E16").Value = "  +8.48%  "
$ws.Range("D17").Value = "2.337.75"
$ws.Range("E17").Value = "  +5.75%  "
$ws.Range("D18").Value = "43.767.36"
$ws.Range("E18").Value = "  +0.40%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.0000110"
$ws.Range("E19").Value = "  +6.06%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.47"
$ws.Range("E20").Value = "  +7.44%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "72.77"
$ws.Range("E21").Value = "  +3.30%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "239.46"
$ws.Range("E22").Value = "  +2.71%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.28"
$ws.Range("E23").Value = "  -2.21%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.44"
$ws.Range("E24").Value = "  +6.26%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.54"
$ws.Range("E26").Value = "  +1.20%  "
$ws.Range("E27").Value = "  +1.99%  "
$ws.Range("E28").Value = "  -1.66%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.28"
$ws.Range("E29").Value = "  +0.51%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "22.61"
$ws.Range("E30").Value = "  +8.53%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "38.06"
$ws.Range("E31").Value = "  -5.20%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "173.72"
$ws.Range("E32").Value = "  +0.56%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0903"
$ws.Range("E33").Value = "  -2.91%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.51"
$ws.Range("E34").Value = "  +0.37%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.126"
$ws.Range("E35").Value = "  +2.46%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0361"
$ws.Range("E36").Value = "  +3.01%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.109"
$ws.Range("E37").Value = "  -2.50%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.41"
$ws.Range("E38").Value = "  +2.46%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.41"
$ws.Range("E39").Value = "  -3.96%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.39"
$ws.Range("E40").Value = "  +9.92%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.239"
$ws.Range("E41").Value = "  +10.61%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.40"
$ws.Range("E42").Value = "  +20.97%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "12.16"
$ws.Range("E43").Value = "  -3.73%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "9.27"
$ws.Range("E44").Value = "  +10.58%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "62.50"
$ws.Range("E45").Value = "  -1.05%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "5.38"
$ws.Range("E46").Value = "  +0.52%  "
$ws.Range("E47").Value = "  +5.35%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "100.62"
$ws.Range("E48").Value = "  +0.25%  "
$ws.Range("E49").Value = "  +1.11%  "
$ws.Range("D50").Value = "2.564.50"
$ws.Range("E50").Value = "  +4.30%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.187"
$ws.Range("E51").Value = "  +14.79%  "
